$d = $word.ActiveDocument

$pairs = @(
    @("91-72=19", "61-44=17"),
    @("38+25=63", "39+57=96"),
    @("50-12=38", "83+8=91"),
    @("5+26=31", "55-36=19"),
    @("27+25=52", "73-15=58"),
    @("59+39=98", "6+18=24"),
    @("41-23=18", "56-7=49"),
    @("49+23=72", "96-8=88"),
    @("7+66=73", "26+29=55"),
    @("75-9=66", "4+58=62"),
    @("56-37=19", "37+7=44"),
    @("32-16=16", "56+35=91"),
    @("96-58=38", "29+12=41"),
    @("34-18=16", "84-69=15"),
    @("72-18=54", "45+38=83"),
    @("48+44=92", "91-14=77"),
    @("15+26=41", "13+69=82"),
    @("16+66=82", "81-59=22"),
    @("93-87=6", "60-59=1"),
    @("29+15=44", "85+9=94"),
    @("74-59=15", "92-48=44"),
    @("19+23=42", "55+39=94"),
    @("36+39=75", "71-37=34"),
    @("43-35=8", "26+57=83"),
    @("84-76=8", "75+7=82"),
    @("60-44=16", "61-4=57"),
    @("9+3=12", "93-29=64"),
    @("78+13=91", "9+56=65"),
    @("78+7=85", "90-51=39"),
    @("30-19=11", "35-27=8"),
    @("38+29=67", "16+15=31"),
    @("92-39=53", "90-65=25"),
    @("42-36=6", "60-59=1"),
    @("38+15=53", "90-19=71"),
    @("26-8=18", "49+12=61"),
    @("62-45=17", "35+7=42"),
    @("82-25=57", "38+28=66"),
    @("70-42=28", "30-28=2"),
    @("55+7=62", "90-86=4"),
    @("39+23=62", "74-7=67"),
    @("74+8=82", "84-75=9"),
    @("83-55=28", "87-68=19"),
    @("54+38=92", "37+46=83"),
    @("56+36=92", "83-79=4"),
    @("58-39=19", "80-5=75"),
    @("16+38=54", "17+58=75"),
    @("53-47=6", "53-35=18"),
    @("50-48=2", "81-32=49"),
    @("15+59=74", "87-49=38"),
    @("35-6=29", "57+27=84"),
    @("85-67=18", "44-36=8"),
    @("8+48=56", "84-58=26"),
    @("91-15=76", "52-15=37"),
    @("26+56=82", "70-29=41"),
    @("61-36=25", "94-89=5"),
    @("70-15=55", "69+15=84"),
    @("62+9=71", "74-8=66"),
    @("61-3=58", "38+8=46"),
    @("51-25=26", "12+59=71"),
    @("7+17=24", "66+7=73"),
    @("44+17=61", "82-58=24"),
    @("79+9=88", "41-36=5"),
    @("68+24=92", "46-29=17"),
    @("21-17=4", "71-52=19"),
    @("81-22=59", "25+7=32"),
    @("83-69=14", "41-3=38"),
    @("41-38=3", "92-65=27"),
    @("17+6=23", "47+38=85"),
    @("93-54=39", "37+25=62"),
    @("37+15=52", "22+9=31"),
    @("41-13=28", "43+29=72"),
    @("26+69=95", "36+27=63"),
    @("46+17=63", "80-2=78"),
    @("51-18=33", "94-46=48"),
    @("9+62=71", "8+63=71"),
    @("65+27=92", "42+39=81"),
    @("52-43=9", "41-6=35"),
    @("9+13=22", "4+77=81"),
    @("37+19=56", "84-79=5"),
    @("95-67=28", "95-78=17"),
    @("63+18=81", "50-18=32"),
    @("65-16=49", "75-16=59"),
    @("17+28=45", "6+27=33"),
    @("25+38=63", "85-49=36"),
    @("94-45=49", "14-9=5"),
    @("18+43=61", "61-23=38"),
    @("93-44=49", "19+74=93"),
    @("50-3=47", "83+8=91"),
    @("57+5=62", "47+19=66"),
    @("70-61=9", "19+12=31"),
    @("19+2=21", "59+8=67"),
    @("50-16=34", "75-6=69"),
    @("65+7=72", "80-44=36"),
    @("55+19=74", "70-28=42"),
    @("56+17=73", "31-25=6"),
    @("97-18=79", "38+38=76"),
    @("85-59=26", "87-8=79"),
    @("31-14=17", "68-59=9"),
    @("59+22=81", "94-6=88"),
    @("62-7=55", "83-56=27")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
